$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.994.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.741.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.77%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4993"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +6.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3601"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.65%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07268"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.061"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.28"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.968"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.740.96"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.854"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001036"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06374"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.05%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.722"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.068.10"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.045"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.46"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.96"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.939.47"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.148"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.67"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.050"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09551"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.574"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.391"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02205"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05877"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.04"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.434"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.88%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.766"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6025"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.111"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.530"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.80"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.601"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5649"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.863"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.104"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06668"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.90%  "

